$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D to fit the new, longer "Cause" text that is about to be added.
$ws.Columns.Item(4).ColumnWidth = 46.8

# --- New shared strings must be introduced in the same order they first appear
# in the final sharedStrings table (indices 50-53), so set these cells first,
# in that exact order: B16, D16, E15, E16.
$ws.Range("B16").Value = "player can still get cryo core even after solving mastermind puzzle"
$ws.Range("D16").Value = "combat logic doesn" + [char]0x2019 + "t check for cryocore in inventory"
$ws.Range("E15").Value = "added enemy death check to fight command"
$ws.Range("E16").Value = "added cryocore check in combatsystem"

# --- Remaining cells reuse already-existing shared strings.
$ws.Range("A16").Value = 15
$ws.Range("C16").Value = "ricky"
$ws.Range("F15").Value = "fixed"
$ws.Range("F16").Value = "fixed"

# --- View state: select E15 as the active cell.
$ws.Range("E15").Select()
